$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.435.99"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.885.04"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  -0.80%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -4.91%  "
$ws.Range("E7").Value = "  -0.89%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.05"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.351"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.79%  "
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("E11").Value = "  -2.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "13.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "2.158.65"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.739"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  -1.53%  "
$ws.Range("D16").Value = "1.875.48"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("D17").Value = "35.388.52"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "73.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "0.0₃0821"
$ws.Range("E19").Value = "  -3.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "244.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("E22").Value = "  -4.55%  "
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("E24").Value = "  +3.45%  "
$ws.Range("E25").Value = "  -9.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.60%  "
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("E29").Value = "  -4.14%  "
$ws.Range("D30").Value = "4.128.47"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("E32").Value = "  -3.05%  "
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.90%  "
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.847"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -16.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0684"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "97.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.77%  "
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("E43").Value = "  -4.35%  "
$ws.Range("D44").Value = "1.292.74"
$ws.Range("E44").Value = "  -3.84%  "
$ws.Range("E45").Value = "  -5.11%  "
$ws.Range("E46").Value = "  +6.91%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.98"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.82%  "
